# Applies the template.xlsx edit described in the commit:
#  - trims the trailing space off the "area1 " sheet tab name
#  - moves the live selection on a few sheets (venue, venue-layout, area2)
#  - switches the active/selected tab from "venue" to "welcome"

$wb = $excel.ActiveWorkbook

# 1. Rename "area1 " -> "area1" (drop the stray trailing space in the tab name).
$areaOne = $wb.Worksheets.Item("area1 ")
$areaOne.Name = "area1"

# 2. Move the saved selection on the sheets whose activeCell/sqref changed.
#    Selecting a range also makes that sheet the active one while the call
#    runs, so we do the non-final sheets first and finish on "welcome" so it
#    ends up as the workbook's active/selected tab (matching activeTab="0",
#    the implicit default once "welcome" - the first sheet - is selected).
$wb.Worksheets.Item("venue").Range("C9").Select()
$wb.Worksheets.Item("venue-layout").Range("I7").Select()
$wb.Worksheets.Item("area2").Range("J21").Select()
$wb.Worksheets.Item("welcome").Range("E27").Select()
